# Remove the trailing "Requisitos" section (its Heading2 title paragraph
# plus the following ListBullet paragraph naming the prerequisite course).
$d = $word.ActiveDocument

$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$headingPar = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Requisitos") {
        $headingPar = $p
        break
    }
}

if ($headingPar -ne $null) {
    $deleteRange = $d.Range($headingPar.Range.Start, $d.Content.End)
    $deleteRange.Delete()
}
